$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header cell C1 from "Voice" to "Voice ID"
$ws.Range("C1").Value = "Voice ID"

# Move selection to E4 (matches the saved selection state in the diff)
$ws.Range("E4").Select()
